$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: rename header row suffixes (_old -> _FV2304, _new -> _FV2310)
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2304")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2310")
}

# Step 2: create Excel Table over A1:U91 (no particular table style)
$range = $ws.Range("A1:U91")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# Step 3: freeze top row (row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
